# Add three new worksheets (SpecsDataCalib1/2/3), each a copy of the
# existing "SpecsDataCalib" sheet, but with cells AA2/AB2/AC2 updated to
# new calibration values.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("SpecsDataCalib")

$newNames = @("SpecsDataCalib1", "SpecsDataCalib2", "SpecsDataCalib3")
$values = @(
    @{ AA2 = 0.1283181294047265; AB2 = 1.338921625438018;  AC2 = 0.02335669498174712 },
    @{ AA2 = 0.1283181294047265; AB2 = 1.453444456925703;  AC2 = 0.320394447286111   },
    @{ AA2 = 0.1194181269674673; AB2 = 1.346468636883841;  AC2 = 0.3163577515461272  }
)

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = $newNames[$i]

    $v = $values[$i]
    $newSheet.Range("AA2").Value = $v.AA2
    $newSheet.Range("AB2").Value = $v.AB2
    $newSheet.Range("AC2").Value = $v.AC2
}
